$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 254.18182
$ws.Range("J6").Value = 298.57144
$ws.Range("L6").Value = 895.71432
$ws.Range("N6").Value = -1119.71432
$ws.Range("H40").Value = 149785.58
$ws.Range("J40").Value = 174333.17
$ws.Range("L40").Value = 174333.17
$ws.Range("N40").Value = -174683.17
$ws.Range("H57").Value = 49121.668
$ws.Range("J57").Value = 49121.668
$ws.Range("L57").Value = 147365.004
$ws.Range("N57").Value = -148363.004
$ws.Range("H100").Value = 2187.9375
$ws.Range("J100").Value = 3544.6667
$ws.Range("L100").Value = 3544.6667
$ws.Range("N100").Value = -4626.6667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2432.25
$ws.Range("J2").Value = 3465
$ws.Range("L2").Value = 3465
$ws.Range("N2").Value = -3691
$ws.Range("H32").Value = 3852596.2
$ws.Range("I32").Value = 4086336.8
$ws.Range("K32").Value = 4086336.8
$ws.Range("M32").Value = -4086049.8
$ws.Range("H45").Value = 2404.9412
$ws.Range("I45").Value = 1740.125
$ws.Range("K45").Value = 1740.125
$ws.Range("M45").Value = -1363.125
$ws.Range("H51").Value = 100000
$ws.Range("J51").Value = 100000
$ws.Range("L51").Value = 100000
$ws.Range("N51").Value = -101512
$ws.Range("H61").Value = 8957
$ws.Range("I61").Value = 2300.875
$ws.Range("K61").Value = 2300.875
$ws.Range("M61").Value = -2088.875
$ws.Range("H110").Value = 30304380
$ws.Range("J110").Value = 66668036
$ws.Range("L110").Value = 66668036
$ws.Range("N110").Value = -66672126
$ws.Range("H116").Value = 2432.25
$ws.Range("J116").Value = 3465
$ws.Range("L116").Value = 3465
$ws.Range("N116").Value = -8053
$ws.Range("H132").Value = 5043.014
$ws.Range("I132").Value = 3535.4807
$ws.Range("J132").Value = 9168.895
$ws.Range("K132").Value = 10606.4421
$ws.Range("L132").Value = 27506.685
$ws.Range("M132").Value = -8076.4421
$ws.Range("N132").Value = -32566.685
$ws.Range("H136").Value = 8957
$ws.Range("I136").Value = 2300.875
$ws.Range("K136").Value = 6902.625
$ws.Range("M136").Value = -4352.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2432.25
$ws.Range("J3").Value = 3465
$ws.Range("L3").Value = 3465
$ws.Range("N3").Value = -3693
$ws.Range("H20").Value = 166666670
$ws.Range("I20").Value = 166666670
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 166666670
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -166666423
$ws.Range("N20").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6786.98
$ws.Range("J31").Value = 9243.727999999999
$ws.Range("L31").Value = 9243.727999999999
$ws.Range("N31").Value = -9833.727999999999
$ws.Range("H34").Value = 6786.98
$ws.Range("J34").Value = 9243.727999999999
$ws.Range("L34").Value = 9243.727999999999
$ws.Range("N34").Value = -9647.727999999999
$ws.Range("H99").Value = 6353.25
$ws.Range("J99").Value = 7802.8
$ws.Range("L99").Value = 7802.8
$ws.Range("N99").Value = -10798.8
$ws.Range("H107").Value = 1205.4482
$ws.Range("I107").Value = 286.54544
$ws.Range("K107").Value = 286.54544
$ws.Range("M107").Value = 1633.45456
$ws.Range("H126").Value = 6353.25
$ws.Range("J126").Value = 7802.8
$ws.Range("L126").Value = 23408.4
$ws.Range("N126").Value = -28348.4
$ws.Range("H134").Value = 5987.241
$ws.Range("I134").Value = 1494
$ws.Range("K134").Value = 4482
$ws.Range("M134").Value = -1947
$ws.Range("H139").Value = 80741.71000000001
$ws.Range("I139").Value = 19000
$ws.Range("J139").Value = 91032
$ws.Range("K139").Value = 19000
$ws.Range("L139").Value = 91032
$ws.Range("M139").Value = -13860
$ws.Range("N139").Value = -101312
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 78039.71000000001
$ws.Range("I139").Value = 146079.86
$ws.Range("J139").Value = 9999.571
$ws.Range("K139").Value = 438239.58
$ws.Range("L139").Value = 29998.713
$ws.Range("M139").Value = -433099.58
$ws.Range("N139").Value = -40278.713
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1666751.1
$ws.Range("J2").Value = 5000052
$ws.Range("L2").Value = 5000052
$ws.Range("N2").Value = -5000278
$ws.Range("H70").Value = 333341340
$ws.Range("I70").Value = 1000000000
$ws.Range("J70").Value = 12000
$ws.Range("K70").Value = 1000000000
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -999999730
$ws.Range("N70").Value = -12540
$ws.Range("H73").Value = 333341340
$ws.Range("I73").Value = 1000000000
$ws.Range("J73").Value = 12000
$ws.Range("K73").Value = 1000000000
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -999999064
$ws.Range("N73").Value = -13872
$ws.Range("H122").Value = 46818.43
$ws.Range("I122").Value = 55957.348
$ws.Range("K122").Value = 167872.044
$ws.Range("M122").Value = -165422.044
$ws.Range("H132").Value = 4031.5
$ws.Range("I132").Value = 2013.0416
$ws.Range("J132").Value = 10086.875
$ws.Range("K132").Value = 6039.1248
$ws.Range("L132").Value = 30260.625
$ws.Range("M132").Value = -3509.1248
$ws.Range("N132").Value = -35320.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4753.1763
$ws.Range("I7").Value = 2717.111
$ws.Range("J7").Value = 7043.75
$ws.Range("K7").Value = 2717.111
$ws.Range("L7").Value = 7043.75
$ws.Range("M7").Value = -2605.111
$ws.Range("N7").Value = -7267.75
$ws.Range("H22").Value = 1487.8
$ws.Range("I22").Value = 420.22223
$ws.Range("J22").Value = 3089.1667
$ws.Range("K22").Value = 420.22223
$ws.Range("L22").Value = 3089.1667
$ws.Range("M22").Value = -125.22223
$ws.Range("N22").Value = -3679.1667
$ws.Range("H27").Value = 1487.8
$ws.Range("I27").Value = 420.22223
$ws.Range("J27").Value = 3089.1667
$ws.Range("K27").Value = 420.22223
$ws.Range("L27").Value = 3089.1667
$ws.Range("M27").Value = -313.22223
$ws.Range("N27").Value = -3303.1667
$ws.Range("H42").Value = 8000
$ws.Range("I42").Value = 8000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -7437
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 5000
$ws.Range("N45").Value = -5814
$ws.Range("H46").Value = 2215.805
$ws.Range("I46").Value = 1678.9678
$ws.Range("K46").Value = 1678.9678
$ws.Range("M46").Value = -1490.9678
$ws.Range("H49").Value = 8000
$ws.Range("I49").Value = 8000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 8000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -7853
$ws.Range("N49").ClearContents()
$ws.Range("H100").Value = 3994.476
$ws.Range("I100").Value = 2911
$ws.Range("K100").Value = 2911
$ws.Range("M100").Value = -2370
$ws.Range("H122").Value = 3445.2415
$ws.Range("I122").Value = 1921.4
$ws.Range("K122").Value = 5764.200000000001
$ws.Range("M122").Value = -3314.200000000001
$ws.Range("H126").Value = 4753.1763
$ws.Range("I126").Value = 2717.111
$ws.Range("J126").Value = 7043.75
$ws.Range("K126").Value = 8151.333
$ws.Range("L126").Value = 21131.25
$ws.Range("M126").Value = -5681.333
$ws.Range("N126").Value = -26071.25
$ws.Range("H132").Value = 6765.4883
$ws.Range("I132").Value = 3550.889
$ws.Range("K132").Value = 10652.667
$ws.Range("M132").Value = -8122.667000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 133338850
$ws.Range("I62").Value = 142861870
$ws.Range("K62").Value = 142861870
$ws.Range("M62").Value = -142861246
$ws.Range("H65").Value = 133338850
$ws.Range("I65").Value = 142861870
$ws.Range("K65").Value = 714309350
$ws.Range("M65").Value = -714306230
$ws.Range("H126").Value = 1248.2222
$ws.Range("I126").Value = 1286.8
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 3860.4
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -1390.4
$ws.Range("N126").Value = -8540
